$d = $word.ActiveDocument

# Locate the target paragraph robustly (don't hardcode an index).
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "By investigating how we can create insurance packages*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Output "Target paragraph not found"
} else {
    # Range covering the paragraph's content, excluding the trailing
    # paragraph-mark character, so InsertXML replaces the runs in place
    # (rather than splicing in a whole new paragraph).
    $full = $target.Range
    $body = $d.Range($full.Start, $full.End - 1)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r><w:t>By investigating how we can create insurance packages that are more affordable</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">educating the masses on the importance having insurance and by researching the market for insights as to who is more likely to buy the product, the relevant </w:t></w:r>' +
        '<w:r><w:t>stakeholders</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">can work toward a well-developed, dynamic and stable insurance industry. </w:t></w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $body.InsertXML($xml)
}
